# Auto-generated: update cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'22.478.68"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = "'1.572.79"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = "'290.37"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('D7').Value = "'0.3698"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.06%  '
$ws.Range('E8').Value = '  +1.53%  '
$ws.Range('D9').Value = "'0.3398"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('D10').Value = "'1.149"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.94%  '
$ws.Range('D11').Value = "'0.07561"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('E13').Value = '  +2.60%  '
$ws.Range('D14').Value = "'6.034"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.91%  '
$ws.Range('D15').Value = "'7.004"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.32%  '
$ws.Range('D16').Value = "'1.572.70"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'0.00001124"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('D18').Value = "'90.56"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('D19').Value = "'0.06788"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').Value = "'1.001"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').Value = "'6.372"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.00%  '
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').Value = "'12.18"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.27%  '
$ws.Range('D24').Value = "'22.481.36"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').Value = "'2.368"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = "'2.652"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.69%  '
$ws.Range('D27').Value = "'20.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = "'149.82"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Value = "'5.059"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('D30').Value = "'124.91"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = "'1.750.72"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').Value = "'1.069"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.20%  '
$ws.Range('D33').Value = "'6.230"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.79%  '
$ws.Range('D34').Value = "'2.017"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('D35').Value = "'9.862"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('D36').Value = "'0.08397"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = "'0.02483"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.40%  '
$ws.Range('D38').Value = "'0.2305"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('D39').Value = "'1.348"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.97%  '
$ws.Range('D40').Value = "'0.06541"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.01%  '
$ws.Range('D41').Value = "'5.443"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('D42').Value = "'11.32"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.69%  '
$ws.Range('D43').Value = "'0.6243"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.24%  '
$ws.Range('D44').Value = "'14.13"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').Value = "'1.001"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').Value = "'3.795"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('D47').Value = "'0.5884"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.37%  '
$ws.Range('D48').Value = "'2.071"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').Value = "'127.60"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.07%  '
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('D51').Value = "'0.07308"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.15%  '
